$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New teacher rows (12-21), columns: A Faculty id, B Password, C first_name,
# D last_name, E email, F gender, G Courses Taught ---
# Data entered column-by-column (matching the order new shared strings were
# introduced in the source workbook): G, then C+D per row, then B, then A, then E,
# and finally F (gender - reuses the existing "Male" string).

# Column G - Courses Taught
$ws.Cells.Item(12, 7).Value = "English Composition"
$ws.Cells.Item(13, 7).Value = "Linear Algebra"
$ws.Cells.Item(14, 7).Value = "Discrete Math"
$ws.Cells.Item(15, 7).Value = "Digital Logic"
$ws.Cells.Item(16, 7).Value = "Operating Systems"
$ws.Cells.Item(17, 7).Value = "Programming Language concepts"
$ws.Cells.Item(18, 7).Value = "Accounting"
$ws.Cells.Item(19, 7).Value = "Object Oriented Programming "
$ws.Cells.Item(20, 7).Value = "Macro Economics"
$ws.Cells.Item(21, 7).Value = "Theory of Computation"

# Columns C & D - first_name / last_name (row by row)
$ws.Cells.Item(12, 3).Value = "Peter"
$ws.Cells.Item(12, 4).Value = "Parker"
$ws.Cells.Item(13, 3).Value = "Bruce"
$ws.Cells.Item(13, 4).Value = "Wayne"
$ws.Cells.Item(14, 3).Value = "Clark"
$ws.Cells.Item(14, 4).Value = "Kent"
$ws.Cells.Item(15, 3).Value = "Jean "
$ws.Cells.Item(15, 4).Value = "Gray"
$ws.Cells.Item(16, 3).Value = "Charles"
$ws.Cells.Item(16, 4).Value = "Xavier"
$ws.Cells.Item(17, 3).Value = "Nauroto"
$ws.Cells.Item(17, 4).Value = "Uzamaki"
$ws.Cells.Item(18, 3).Value = "James"
$ws.Cells.Item(18, 4).Value = "Harden"
$ws.Cells.Item(19, 3).Value = "Mohammad"
$ws.Cells.Item(19, 4).Value = "Ali"
$ws.Cells.Item(20, 3).Value = "Peter"
$ws.Cells.Item(20, 4).Value = "Parker"
$ws.Cells.Item(21, 3).Value = "Charles"
$ws.Cells.Item(21, 4).Value = "Xavier"

# Column B - Password
$ws.Cells.Item(12, 2).Value = "teacher11"
$ws.Cells.Item(13, 2).Value = "teacher12"
$ws.Cells.Item(14, 2).Value = "teacher13"
$ws.Cells.Item(15, 2).Value = "teacher14"
$ws.Cells.Item(16, 2).Value = "teacher15"
$ws.Cells.Item(17, 2).Value = "teacher16"
$ws.Cells.Item(18, 2).Value = "teacher17"
$ws.Cells.Item(19, 2).Value = "teacher18"
$ws.Cells.Item(20, 2).Value = "teacher19"
$ws.Cells.Item(21, 2).Value = "teacher20"

# Column A - Faculty id
$ws.Cells.Item(12, 1).Value = "parkerp2"
$ws.Cells.Item(13, 1).Value = "waynep2"
$ws.Cells.Item(14, 1).Value = "kentp2"
$ws.Cells.Item(15, 1).Value = "grayp2"
$ws.Cells.Item(16, 1).Value = "xavierp2"
$ws.Cells.Item(17, 1).Value = "uzamakip2"
$ws.Cells.Item(18, 1).Value = "hardenp2"
$ws.Cells.Item(19, 1).Value = "alip2"
$ws.Cells.Item(20, 1).Value = "parkerp2"
$ws.Cells.Item(21, 1).Value = "xavierp2"

# Column E - email
$ws.Cells.Item(12, 5).Value = "pparker@google.com.au"
$ws.Cells.Item(13, 5).Value = "bwayne7@yandex.ru"
$ws.Cells.Item(14, 5).Value = "kclark12@free.fr"
$ws.Cells.Item(15, 5).Value = "jgray8@jailbum.net"
$ws.Cells.Item(16, 5).Value = "xcharly19@ihg.com"
$ws.Cells.Item(17, 5).Value = "unauroto20@umich.edu"
$ws.Cells.Item(18, 5).Value = "jharden@rambler.ru"
$ws.Cells.Item(19, 5).Value = "amohamm@jailbum.net"
$ws.Cells.Item(20, 5).Value = "pparker@google.com.au"
$ws.Cells.Item(21, 5).Value = "xcharly19@ihg.com"

# Column F - gender (all Male, reuses existing shared string)
$ws.Cells.Item(12, 6).Value = "Male"
$ws.Cells.Item(13, 6).Value = "Male"
$ws.Cells.Item(14, 6).Value = "Male"
$ws.Cells.Item(15, 6).Value = "Male"
$ws.Cells.Item(16, 6).Value = "Male"
$ws.Cells.Item(17, 6).Value = "Male"
$ws.Cells.Item(18, 6).Value = "Male"
$ws.Cells.Item(19, 6).Value = "Male"
$ws.Cells.Item(20, 6).Value = "Male"
$ws.Cells.Item(21, 6).Value = "Male"

# The "Courses Taught" column now needs to be wider to fit the longer course
# names (user manually widened it, dropping the old auto best-fit).
$ws.Columns.Item(7).ColumnWidth = 29

# Selection left where the user's cursor ended up after entering the new rows.
$ws.Range("E23").Select()
